$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each updated cell as text (apostrophe-prefix forces text entry,
# matching the original inlineStr/text cell type), then reset the cell
# style back to Normal so we do not leave a stray "Text" number-format
# style applied to the cell (the source cells carry no explicit style).
function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell 'D2' '62.923.86'
Set-TextCell 'E2' '  +2.08%  '
Set-TextCell 'D3' '3.467.66'
Set-TextCell 'D4' '0.999'
Set-TextCell 'E4' '  +0.02%  '
Set-TextCell 'D5' '577.15'
Set-TextCell 'E5' '  +0.12%  '
Set-TextCell 'D6' '147.75'
Set-TextCell 'E6' '  +3.59%  '
Set-TextCell 'D7' '3.468.96'
Set-TextCell 'E7' '  +2.17%  '
Set-TextCell 'E8' '  -0.06%  '
Set-TextCell 'E9' '  +1.35%  '
Set-TextCell 'D10' '7.64'
Set-TextCell 'E10' '  +0.17%  '
Set-TextCell 'E11' '  +1.25%  '
Set-TextCell 'E12' '  +4.45%  '
Set-TextCell 'D13' '4.061.15'
Set-TextCell 'E13' '  +2.17%  '
Set-TextCell 'D14' '29.74'
Set-TextCell 'E14' '  +6.42%  '
Set-TextCell 'E15' '  +2.75%  '
Set-TextCell 'D16' '3.474.17'
Set-TextCell 'E16' '  +1.94%  '
Set-TextCell 'E17' '  +0.24%  '
Set-TextCell 'D18' '62.889.96'
Set-TextCell 'E18' '  +1.96%  '
Set-TextCell 'E19' '  +3.32%  '
Set-TextCell 'D20' '14.34'
Set-TextCell 'E20' '  +5.09%  '
Set-TextCell 'D21' '9.20'
Set-TextCell 'E21' '  +1.02%  '
Set-TextCell 'D22' '387.15'
Set-TextCell 'E22' '  -0.49%  '
Set-TextCell 'B23' 'Litecoin'
Set-TextCell 'C23' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 'D23' '74.66'
Set-TextCell 'E23' '  +0.05%  '
Set-TextCell 'B24' 'Polygon'
Set-TextCell 'C24' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 'D24' '0.555'
Set-TextCell 'E24' '  +1.32%  '
Set-TextCell 'E25' '  -0.12%  '
Set-TextCell 'D26' '3.608.21'
Set-TextCell 'E26' '  +2.12%  '
Set-TextCell 'E27' '  +1.19%  '
Set-TextCell 'D28' '0.180'
Set-TextCell 'E28' '  -1.19%  '
Set-TextCell 'D29' '7.55'
Set-TextCell 'E29' '  +2.31%  '
Set-TextCell 'D30' '1.00'
Set-TextCell 'E30' '  +0.32%  '
Set-TextCell 'E31' '  +2.06%  '
Set-TextCell 'E32' '  -1.06%  '
Set-TextCell 'E34' '  -2.79%  '
Set-TextCell 'D35' '23.61'
Set-TextCell 'E35' '  +1.08%  '
Set-TextCell 'D36' '32.02'
Set-TextCell 'E36' '  +18.81%  '
Set-TextCell 'E37' '  +3.01%  '
Set-TextCell 'E38' '  +1.62%  '
Set-TextCell 'D39' '169.87'
Set-TextCell 'E39' '  +0.88%  '
Set-TextCell 'E40' '  +5.37%  '
Set-TextCell 'D41' '3.505.24'
Set-TextCell 'E41' '  +2.27%  '
Set-TextCell 'D42' '0.0751'
Set-TextCell 'E42' '  -1.33%  '
Set-TextCell 'D43' '0.798'
Set-TextCell 'E43' '  +1.93%  '
Set-TextCell 'D44' '42.33'
Set-TextCell 'E44' '  -0.27%  '
Set-TextCell 'E45' '  +0.19%  '
Set-TextCell 'E46' '  +2.48%  '
Set-TextCell 'E47' '  +3.25%  '
Set-TextCell 'D48' '2.608.60'
Set-TextCell 'E48' '  +5.36%  '
Set-TextCell 'D49' '2.25'
Set-TextCell 'E49' '  +10.23%  '
Set-TextCell 'D50' '22.88'
Set-TextCell 'E50' '  +1.07%  '
Set-TextCell 'E51' '  +0.94%  '
